# "easter adjustements to sync to vensim"
# Updates to the initial_species sheet: rename METABOLIC_RATE -> KC_METABOLIC_RATE,
# adjust a handful of input cells (D26, D28, row 29 for columns C:I), highlight the
# TREE_CONTINENTAL column header (D5) in yellow, and flag the corrected D28 cell with
# the workbook's existing "warning" (red-on-light) style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("initial_species")

# --- Rename the METABOLIC_RATE row label to KC_METABOLIC_RATE ---
$ws.Range("A26").Value = "KC_METABOLIC_RATE"

# --- Adjust KC_METABOLIC_RATE (row26) for TREE_CONTINENTAL (col D) ---
$ws.Range("D26").Value = 0.09

# --- Adjust TEMPERATURE_RANGE (row28) for TREE_CONTINENTAL (col D) ---
$ws.Range("D28").Value = 30

# Flag the corrected cell using the same "warning" style already used on
# row 27/28 for the secondary columns (red font on the light accent fill).
$warnSource = $ws.Range("J28")
$warnSource.Copy()
$ws.Range("D28").PasteSpecial(-4122)

# --- Adjust ANABOLISM_BIOMASS_PER_WATER_L (row29) for columns C..I ---
$ws.Range("C29").Value = 0.04
$ws.Range("D29").Value = 0.01
$ws.Range("E29").Value = 0.06
$ws.Range("F29").Value = 0.02
$ws.Range("G29").Value = 0.04
$ws.Range("H29").Value = 0.04
$ws.Range("I29").Value = 0.04

# --- Highlight the TREE_CONTINENTAL header (D5) in yellow ---
$ws.Range("D5").Interior.Color = 65535

# --- Mirror the final on-screen selection left by the edit ---
$ws.Range("D28").Select()

$excel.CalculateFull()
